$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Français" column (B) held one grade per row (B2:B7). The user now
# wants a single averaged/converted note, kept on row 2, with the other
# rows freed up so the user can pick where to put data next.
$ws.Range("B2").Value = 14

# Free up B3:B7 (the grades that used to live there move/condense into B2).
$ws.Range("B3:B7").ClearContents()

# Move the active selection back to A1 (top of sheet).
$ws.Range("A1").Select()
